$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93 (existing row 93 and below shift down to 94, etc.)
$ws.Rows.Item(93).Insert()

# The new row 93 gets a copy of the data that row 92 held before this edit
# (same product/market metadata, with the price-point that used to be "latest")
$ws.Cells.Item(93,1).Value = 10
$ws.Cells.Item(93,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(93,3).Value = "La Araucanía"
$ws.Cells.Item(93,4).Value = 44999
$ws.Cells.Item(93,5).Value = 9
$ws.Cells.Item(93,6).Value = "Fruta"
$ws.Cells.Item(93,7).Value = 100108
$ws.Cells.Item(93,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(93,9).Value = 100108003
$ws.Cells.Item(93,10).Value = "Maracuyá"
$ws.Cells.Item(93,11).Value = "Sin especificar"
$ws.Cells.Item(93,12).Value = "Primera"
$ws.Cells.Item(93,13).Value = 8
$ws.Cells.Item(93,14).Value = 60000
$ws.Cells.Item(93,15).Value = 60000
$ws.Cells.Item(93,16).Value = 60000
$ws.Cells.Item(93,17).Value = "$/caja 18 kilos"
$ws.Cells.Item(93,18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(93,19).Value = 3333
$ws.Cells.Item(93,20).Value = 18

# Give row 93's date cell the same date number format as the rest of column D
$ws.Cells.Item(93,4).NumberFormat = $ws.Cells.Item(92,4).NumberFormat

# Update row 92 with the new, latest weekly price record
$ws.Cells.Item(92,4).Value = 45075
$ws.Cells.Item(92,13).Value = 50
$ws.Cells.Item(92,14).Value = 50000
$ws.Cells.Item(92,15).Value = 50000
$ws.Cells.Item(92,16).Value = 50000
$ws.Cells.Item(92,19).Value = 2778
